$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" column (C) for rows 2-11 from 45174 to 45175
$ws.Range("C2:C11").Value = 45175
